$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A81").Value = "'2025/10/09"
$ws.Range("A81").ClearFormats()
$ws.Range("B81").Value = "木"
$ws.Range("C81").Value = 1
$ws.Range("D81").Value = 14
